$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reset the data block so stale shared-string entries are fully
# dereferenced before writing the refreshed correlation table back in.
$data = $ws.Range("B2:G9")
$data.NumberFormat = "@"
$data.Value = "zzz"
$data.ClearContents()
$data.Style = "Normal"

# --- Refreshed correlation values (updated sample / re-estimated density),
# written column by column, top to bottom.
$cols = [ordered]@{
    B = @("-0.02",   "-0.05",   "-0.01",   "nan", "0.02",  "0.01",  "0.02",  "-0.08")
    C = @("-0.15",   "-0.12",   "-0.01",   "nan", "-0.07", "-0.04", "-0.0",  "-0.22*")
    D = @("-0.24**", "-0.25**", "-0.06",   "nan", "-0.09", "-0.14", "0.07",  "0.0")
    E = @("-0.18",   "-0.17",   "-0.03",   "nan", "-0.11", "-0.14", "-0.07", "-0.14")
    F = @("-0.05",   "-0.04",   "-0.14",   "nan", "-0.16", "-0.16", "-0.16", "0.12")
    G = @("-0.1",    "-0.12",   "-0.26**", "nan", "-0.12", "-0.12", "-0.18", "0.03")
}

foreach ($col in $cols.Keys) {
    $vals = $cols[$col]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $row = 2 + $i
        $cell = $ws.Range("$col$row")
        $cell.NumberFormat = "@"
        $cell.Value = $vals[$i]
        $cell.Style = "Normal"
    }
}
